# The workbook was re-saved a day later (17072025 -> 18072025 folder), the
# T2 value was refreshed, and the active selection moved down to T3.
# Window-chrome / revision GUID bookkeeping (absPath, revisionPtr, xWindow)
# is regenerated internally by Excel on save and isn't something exposed
# through the object model, so we only script the user-visible data change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the daily figure in T2
$ws.Range("T2").Value = 467290

# Move the active selection from T2 to T3, matching the saved cursor position
$ws.Range("T3").Select()
